$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Table1" list object and its calculated "المقابل النقدى" (cash
# equivalent) column, which used to hold the formula
# Table1[[#This Row],[الفئة]]*Table1[[#This Row],[عدد ايام الورادى]]
$lo = $ws.ListObjects.Item(1)
$col5 = $lo.ListColumns.Item(5)

# Correct the attendance-days figure for the first employee row.
$ws.Range("C2").Value = 10

# Convert the calculated column from live formulas to plain static
# values (matching the new, corrected numbers) - row 2 changes from
# 825 (=75*11) to 750 (=75*10); rows 3-4 keep their existing values.
$ws.Range("E2").Value = 750
$ws.Range("E3").Value = 150
$ws.Range("E4").Value = 225

# Apply a "0.00" number format to the whole column (header included),
# which is now no longer driven by the table's calculated formula.
$ws.Range("E1:E4").NumberFormat = "0.00"

# The calculated column formula is no longer valid once the column
# was converted to static values.
try {
    $col5.CalculatedColumnFormula = "#REF!"
} catch {
}

# Reflect the last user selection recorded in the sheet view (entire
# column G was selected).
$ws.Columns("G").Select()
